$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3248.375
$ws.Range("I32").Value = 3566.3333
$ws.Range("K32").Value = 3566.3333
$ws.Range("M32").Value = -3240.3333
$ws.Range("H40").Value = 3999.36
$ws.Range("I40").Value = 5054.4614
$ws.Range("K40").Value = 5054.4614
$ws.Range("M40").Value = -4879.4614
$ws.Range("H51").Value = 3297.5
$ws.Range("J51").Value = 3396.6667
$ws.Range("L51").Value = 3396.6667
$ws.Range("N51").Value = -4364.6667
$ws.Range("H75").Value = 100000
$ws.Range("J75").Value = 100000
$ws.Range("L75").Value = 100000
$ws.Range("N75").Value = -101872
$ws.Range("H78").Value = 100000
$ws.Range("J78").Value = 100000
$ws.Range("L78").Value = 300000
$ws.Range("N78").Value = -309360
$ws.Range("H113").Value = 4043.1428
$ws.Range("I113").Value = 3268.3333
$ws.Range("J113").Value = 4624.25
$ws.Range("K113").Value = 3268.3333
$ws.Range("L113").Value = 4624.25
$ws.Range("M113").Value = -14.33329999999978
$ws.Range("N113").Value = -11132.25
$ws.Range("H116").Value = 348164.34
$ws.Range("J116").Value = 715692.6
$ws.Range("L116").Value = 715692.6
$ws.Range("N116").Value = -722576.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3704668.8
$ws.Range("I61").Value = 3704668.8
$ws.Range("K61").Value = 3704668.8
$ws.Range("M61").Value = -3704456.8
$ws.Range("H74").Value = 6241
$ws.Range("I74").Value = 1973.25
$ws.Range("K74").Value = 1973.25
$ws.Range("M74").Value = -1099.25
$ws.Range("H77").Value = 6241
$ws.Range("I77").Value = 1973.25
$ws.Range("K77").Value = 9866.25
$ws.Range("M77").Value = -5498.25
$ws.Range("H122").Value = 1046.7142
$ws.Range("I122").Value = 973.1667
$ws.Range("K122").Value = 2919.5001
$ws.Range("M122").Value = -469.5001000000002
$ws.Range("H132").Value = 490127.38
$ws.Range("I132").Value = 557506.4399999999
$ws.Range("K132").Value = 1672519.32
$ws.Range("M132").Value = -1669989.32
$ws.Range("H136").Value = 3704668.8
$ws.Range("I136").Value = 3704668.8
$ws.Range("K136").Value = 11114006.4
$ws.Range("M136").Value = -11111456.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 10000
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 1015975.8
$ws.Range("I134").Value = 954634.6
$ws.Range("K134").Value = 2863903.8
$ws.Range("M134").Value = -2861368.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 5633.3335
$ws.Range("I36").Value = 5633.3335
$ws.Range("K36").Value = 5633.3335
$ws.Range("M36").Value = -5245.3335
$ws.Range("H40").Value = 5633.3335
$ws.Range("I40").Value = 5633.3335
$ws.Range("K40").Value = 5633.3335
$ws.Range("M40").Value = -5473.3335
$ws.Range("H58").Value = 1241135.6
$ws.Range("I58").Value = 3090489
$ws.Range("J58").Value = 8233.333000000001
$ws.Range("K58").Value = 3090489
$ws.Range("L58").Value = 8233.333000000001
$ws.Range("M58").Value = -3090286
$ws.Range("N58").Value = -8639.333000000001
$ws.Range("H62").Value = 3132.7778
$ws.Range("I62").Value = 2498.75
$ws.Range("J62").Value = 3640
$ws.Range("K62").Value = 2498.75
$ws.Range("L62").Value = 3640
$ws.Range("M62").Value = -1874.75
$ws.Range("N62").Value = -4888
$ws.Range("H65").Value = 3132.7778
$ws.Range("I65").Value = 2498.75
$ws.Range("J65").Value = 3640
$ws.Range("K65").Value = 12493.75
$ws.Range("L65").Value = 18200
$ws.Range("M65").Value = -9373.75
$ws.Range("N65").Value = -24440
$ws.Range("H134").Value = 5004.1465
$ws.Range("I134").Value = 5309.1387
$ws.Range("K134").Value = 15927.4161
$ws.Range("M134").Value = -13392.4161
$ws.Range("H136").Value = 1241135.6
$ws.Range("I136").Value = 3090489
$ws.Range("J136").Value = 8233.333000000001
$ws.Range("K136").Value = 9271467
$ws.Range("L136").Value = 24699.999
$ws.Range("M136").Value = -9268917
$ws.Range("N136").Value = -29799.999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 273.13333
$ws.Range("I8").Value = 273.13333
$ws.Range("K8").Value = 819.39999
$ws.Range("M8").Value = -680.39999
$ws.Range("H75").Value = 5315.636
$ws.Range("I75").Value = 1749.5
$ws.Range("J75").Value = 6108.1113
$ws.Range("K75").Value = 5248.5
$ws.Range("L75").Value = 18324.3339
$ws.Range("M75").Value = -4250.5
$ws.Range("N75").Value = -20320.3339
$ws.Range("H78").Value = 5315.636
$ws.Range("I78").Value = 1749.5
$ws.Range("J78").Value = 6108.1113
$ws.Range("K78").Value = 15745.5
$ws.Range("L78").Value = 54973.00169999999
$ws.Range("M78").Value = -10753.5
$ws.Range("N78").Value = -64957.00169999999
$ws.Range("H126").Value = 2688.6667
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 2533
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 7599
$ws.Range("M126").Value = -4060
$ws.Range("N126").Value = -17479
$ws.Range("H131").Value = 21462.092
$ws.Range("I131").Value = 954.8333
$ws.Range("J131").Value = 46070.8
$ws.Range("K131").Value = 2864.4999
$ws.Range("L131").Value = 138212.4
$ws.Range("M131").Value = 2175.5001
$ws.Range("N131").Value = -148292.4
$ws.Range("H132").Value = 1334.2858
$ws.Range("I132").Value = 795
$ws.Range("J132").Value = 1550
$ws.Range("K132").Value = 7155
$ws.Range("L132").Value = 13950
$ws.Range("M132").Value = -4625
$ws.Range("N132").Value = -19010

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 75773.13
$ws.Range("I122").Value = 116010.89
$ws.Range("K122").Value = 348032.67
$ws.Range("M122").Value = -345582.67
$ws.Range("H132").Value = 25305904
$ws.Range("I132").Value = 30671950
$ws.Range("K132").Value = 92015850
$ws.Range("M132").Value = -92013320

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1944.1111
$ws.Range("J22").Value = 2416.6667
$ws.Range("L22").Value = 2416.6667
$ws.Range("N22").Value = -3006.6667
$ws.Range("H27").Value = 1944.1111
$ws.Range("J27").Value = 2416.6667
$ws.Range("L27").Value = 2416.6667
$ws.Range("N27").Value = -2630.6667
$ws.Range("H68").Value = 4334.8335
$ws.Range("I68").Value = 6500
$ws.Range("J68").Value = 3252.25
$ws.Range("K68").Value = 6500
$ws.Range("L68").Value = 3252.25
$ws.Range("M68").Value = -5751
$ws.Range("N68").Value = -4750.25
$ws.Range("H71").Value = 4334.8335
$ws.Range("I71").Value = 6500
$ws.Range("J71").Value = 3252.25
$ws.Range("K71").Value = 32500
$ws.Range("L71").Value = 16261.25
$ws.Range("M71").Value = -28756
$ws.Range("N71").Value = -23749.25
$ws.Range("H136").Value = 2573.5652
$ws.Range("I136").Value = 1837.375
$ws.Range("K136").Value = 5512.125
$ws.Range("M136").Value = -2962.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 5000
$ws.Range("I32").Value = 5000
$ws.Range("K32").Value = 5000
$ws.Range("M32").Value = -4683
$ws.Range("H76").Value = 15000
$ws.Range("I76").Value = 15000
$ws.Range("K76").Value = 15000
$ws.Range("M76").Value = -14685
$ws.Range("H79").Value = 15000
$ws.Range("I79").Value = 15000
$ws.Range("K79").Value = 15000
$ws.Range("M79").Value = -13908
$ws.Range("H132").Value = 4027498.8
$ws.Range("I132").Value = 5033186
$ws.Range("K132").Value = 15099558
$ws.Range("M132").Value = -15097028
$ws.Range("H136").Value = 8787.031999999999
$ws.Range("I136").Value = 9824.208000000001
$ws.Range("J136").Value = 5231
$ws.Range("K136").Value = 29472.624
$ws.Range("L136").Value = 15693
$ws.Range("M136").Value = -26922.624
$ws.Range("N136").Value = -20793
